$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data rows (15-20) appended to the Hover data set
$data = @(
    @(14, "2018-10-07 13-59-17.csv", 3, "train", "2018-10-07T14:00:08.000", 6),
    @(15, "2018-10-07 13-59-17.csv", 3, "train", "2018-10-07T14:00:18.000", 6),
    @(16, "2018-10-07 13-59-17.csv", 3, "test",  "2018-10-07T14:00:28.000", 6),
    @(17, "2018-10-07 13-59-17.csv", 3, "train", "2018-10-07T14:00:39.500", 6),
    @(18, "2018-10-07 13-59-17.csv", 3, "test",  "2018-10-07T14:00:59.00",  6),
    @(19, "2018-10-07 13-59-17.csv", 3, "train", "2018-10-07T14:01:22.00",  6)
)

$row = 15
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $ws.Cells.Item($row, 6).Value = $entry[5]
    $row++
}

# Update the selected cell to match the post-edit state
$null = $ws.Range("H15").Select()

Write-Output "done"
